# Insert a new "ObjectName" worksheet before "Sheet1" (so it becomes the
# 3rd sheet, right after BankDetail), populate it with the new report
# headers/values, size its columns, and make it the active/selected sheet
# -- matching the author's "End to End Test suite" commit.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Add($sheet1)
$ws.Name = "ObjectName"

# Header row
$ws.Range("A1").Value = "Company"
$ws.Range("B1").Value = "GLAccount"
$ws.Range("C1").Value = "Vendor"
$ws.Range("D1").Value = "Fiscalyear"
$ws.Range("E1").Value = "BankName"

# Data row
$ws.Range("A2").Value = "Mountain Manufacturing (100)"
$ws.Range("B2").Value = "6655 (BC)"
$ws.Range("C2").Value = "SB-24808 (1103)"
$ws.Range("D2").Value = 2022
$ws.Range("E2").Value = "Bank of Dad"

# Column sizing to fit the new content
$ws.Columns.Item(1).ColumnWidth = 27.5
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 8.833333333333334
$ws.Columns.Item(5).ColumnWidth = 10.666666666666666

# Make the new sheet active with its own selection, like in the source file
$ws.Activate()
$ws.Range("E11").Select()
